$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the auto-generated DSR test data rows (rows 2-4) with new
# randomized values produced by the latest headless test run.

# Column E holds a phone number that is all digits with a leading zero,
# so format it as Text first to keep Excel from coercing it to a number
# and stripping the leading zero.
$ws.Range("E2:E4").NumberFormat = "@"

$ws.Range("C2").Value = "AUTODSR_7A547"
$ws.Range("D2").Value = "Father_0988"
$ws.Range("E2").Value = "03869370300"
$ws.Range("G2").Value = "EMP73754D"

$ws.Range("C3").Value = "AUTODSR_1F1D4"
$ws.Range("D3").Value = "Father_C864"
$ws.Range("E3").Value = "03870941100"
$ws.Range("G3").Value = "EMP375387"

$ws.Range("C4").Value = "AUTODSR_9CB31"
$ws.Range("D4").Value = "Father_2EC4"
$ws.Range("E4").Value = "03871865900"
$ws.Range("G4").Value = "EMP687693"
